$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("Payton Pritchard", "PG,SG", "Boston Celtics"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Kyshawn George", "SG,SF", "Washington Wizards"),
    @("Kevin Porter Jr.", "PG", "Milwaukee Bucks"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Kyle Filipowski", "PF,C", "Utah Jazz"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Russell Westbrook", "PG,SG", "Denver Nuggets"),
    @("Stephon Castle", "PG,SG", "San Antonio Spurs"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# The table now has one fewer row than before (18 rows total incl. header
# instead of 19), so delete the now-obsolete last row.
$ws.Rows.Item(19).Delete()
